$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Study_Types" value "Interventional" is being renamed to "Clinical"
# across all rows where it appears (the data-staging group rows).
$ws.Range("C2").Value = "Clinical"
$ws.Range("C7").Value = "Clinical"
$ws.Range("C12").Value = "Clinical"

# Update the active selection shown in the saved view (was F8, now C13).
$ws.Range("C13").Select()
